$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B=838.1461181640625;  C=0.6288;              D=0.5997999906539917; E=0.9879000186920166; F=0.1492000073194504}
    @{Row=3;  B=956.6784057617188;  C=0.7653;              D=0.7966;             E=1;                   F=0.1378999948501587}
    @{Row=4;  B=661.7163696289062;  C=0.773;               D=0.7822;             E=0.9682999849319458; F=0.2104000002145767}
    @{Row=5;  B=420.0038146972656;  C=0.443;               D=0.365;              E=0.9375;              F=0.05860000103712082}
    @{Row=6;  B=607.0703735351562;  C=0.4803;              D=0.5135;             E=0.9347000122070312; F=0.09529999643564224}
    @{Row=7;  B=541.568115234375;   C=0.5454;              D=0.5611000061035156; E=0.7748000025749207; F=0.2646999955177307}
    @{Row=8;  B=596.7949829101562;  C=0.5348000000000001;  D=0.5403;             E=0.7748000025749207; F=0.3339000046253204}
    @{Row=9;  B=4621.97802734375;   C=0.5959;              D=0.5775;             E=1;                   F=0.05860000103712082}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
}
